# Locate the paragraph that ends the existing "COPIAR LOS DATOS..." block
# ("UPDATE platos SET id_secundario = concat(...) WHERE id_plato > 1") and
# append the new "ACTUALIZAR O INSERTAR REGISTRO ... OTRA COLUMNA EN OTRA
# TABLA" section (heading + UPDATE/SET/WHERE statement) right after it.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute(
    "UPDATE platos SET id_secundario = concat(id_plato, '-', tipo_plato) WHERE id_plato > 1",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

# Expand the found range to its whole paragraph (including the end-of-
# paragraph mark) and collapse to its end so the new content is inserted
# immediately after it.
$anchor.Expand(4) | Out-Null
$anchor.Collapse(0)

$newContentXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>ACTUALIZAR O INSERTAR</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>REGISTRO</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> DE UNA COLUMNA </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">CON LOS REGISTROS DE </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>OTRA COLUMNA EN OTRA TABLA</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>UPDATE pedidos , platosamostrar</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>SET pedidos.precio_plato = platosamostrar.precio_plato</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>WHERE id_pedido > 0</w:t></w:r></w:p>
"@

$null = $anchor.InsertXML($newContentXml)
